# Generate Report for Handoff
#
# A new handoff was generated for "40fbaa99-41b0-40d0-8959-a1553743c22d.md"
# (row 5 in each sheet). Update the "Latest Handoff Date(time)" values that
# the report generator writes: the Overview roll-up date, and the per-locale
# handoff timestamps on the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-30-17 20:30:54"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-17 20:30:51"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-17 20:30:54"
